$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 and 9: coin identities swap (Dogecoin <-> Cardano) with new price/volume data
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "0.3042"
$ws.Range("E8").Value = "  -3.10%  "

$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.07616"
$ws.Range("E9").Value = "  -3.12%  "

# Remaining rows: price (D) and/or volume (E) updates
$ws.Range("D2").Value = "28.879.37"
$ws.Range("E2").Value = "  -1.84%  "
$ws.Range("D3").Value = "1.827.23"
$ws.Range("E3").Value = "  -2.07%  "
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "241.31"
$ws.Range("E5").Value = "  -0.78%  "
$ws.Range("D6").Value = "0.6903"
$ws.Range("E6").Value = "  -2.33%  "
$ws.Range("D7").Value = "0.9998"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D10").Value = "23.42"
$ws.Range("E10").Value = "  -4.30%  "
$ws.Range("D11").Value = "0.07764"
$ws.Range("E11").Value = "  -3.26%  "
$ws.Range("D12").Value = "1.824.03"
$ws.Range("E12").Value = "  -2.80%  "
$ws.Range("D13").Value = "5.042"
$ws.Range("E13").Value = "  -3.02%  "
$ws.Range("D14").Value = "90.19"
$ws.Range("E14").Value = "  -3.42%  "
$ws.Range("D15").Value = "0.6724"
$ws.Range("E15").Value = "  -3.91%  "
$ws.Range("D16").Value = "6.358"
$ws.Range("E16").Value = "  -1.45%  "
$ws.Range("D17").Value = "0.000008265"
$ws.Range("E17").Value = "  -1.21%  "
$ws.Range("D18").Value = "28.860.26"
$ws.Range("E18").Value = "  -1.93%  "
$ws.Range("D19").Value = "241.72"
$ws.Range("E19").Value = "  -4.16%  "
$ws.Range("D20").Value = "2.075.21"
$ws.Range("E20").Value = "  -2.28%  "
$ws.Range("D21").Value = "12.63"
$ws.Range("E21").Value = "  -3.77%  "
$ws.Range("D22").Value = "0.9996"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").Value = "7.394"
$ws.Range("E23").Value = "  -2.83%  "
$ws.Range("D24").Value = "0.9998"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").Value = "0.1466"
$ws.Range("E25").Value = "  -6.07%  "
$ws.Range("D26").Value = "160.82"
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").Value = "8.731"
$ws.Range("E27").Value = "  -3.04%  "
$ws.Range("D28").Value = "18.15"
$ws.Range("E28").Value = "  -3.09%  "
$ws.Range("D29").Value = "1.534"
$ws.Range("E29").Value = "  +2.40%  "
$ws.Range("D30").Value = "4.192"
$ws.Range("E30").Value = "  -3.03%  "
$ws.Range("D31").Value = "4.125"
$ws.Range("E31").Value = "  -3.66%  "
$ws.Range("D32").Value = "1.192"
$ws.Range("E32").Value = "  -1.67%  "
$ws.Range("D33").Value = "0.05093"
$ws.Range("E33").Value = "  -4.04%  "
$ws.Range("D34").Value = "0.7459"
$ws.Range("E34").Value = "  -0.93%  "
$ws.Range("D35").Value = "1.808"
$ws.Range("E35").Value = "  -4.11%  "
$ws.Range("D36").Value = "1.137"
$ws.Range("E36").Value = "  -2.50%  "
$ws.Range("D37").Value = "2.681"
$ws.Range("E37").Value = "  -1.15%  "
$ws.Range("D38").Value = "0.01836"
$ws.Range("E38").Value = "  -2.18%  "
$ws.Range("D39").Value = "1.201.17"
$ws.Range("E39").Value = "  -5.38%  "
$ws.Range("D40").Value = "2.671"
$ws.Range("E40").Value = "  -2.54%  "
$ws.Range("D41").Value = "0.9275"
$ws.Range("E41").Value = "  +3.19%  "
$ws.Range("D42").Value = "108.10"
$ws.Range("E42").Value = "  -1.01%  "
$ws.Range("D43").Value = "0.9993"
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("D44").Value = "0.5161"
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("D45").Value = "1.974.27"
$ws.Range("E45").Value = "  -2.83%  "
$ws.Range("D46").Value = "9.453"
$ws.Range("E46").Value = "  -1.27%  "
$ws.Range("D47").Value = "0.00000000121"
$ws.Range("E47").Value = "  -5.58%  "
$ws.Range("D48").Value = "5.208"
$ws.Range("E48").Value = "  -12.87%  "
$ws.Range("D49").Value = "1.722"
$ws.Range("E49").Value = "  -3.68%  "
$ws.Range("D50").Value = "61.56"
$ws.Range("E50").Value = "  -13.71%  "
$ws.Range("D51").Value = "0.4166"
$ws.Range("E51").Value = "  -3.28%  "
